$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9826360940933228
$ws.Range("B1").Value = 1.438040614128113
$ws.Range("C1").Value = 3.122458934783936
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.511973738670349
